$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 44, shifting rows 44:67 down to 46:69
# (EntireRow insert preserves formatting of the row below, matching the
# existing style used by column D - the date-formatted style index).
$ws.Rows.Item(44).Insert()
$ws.Rows.Item(44).Insert()

# Populate the two newly inserted rows (44 and 45) with the new records.
# Row 44: Melón, Calameño, Segunda
$ws.Cells.Item(44, 1).Value = 1
$ws.Cells.Item(44, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(44, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(44, 4).Value = 44904
$ws.Cells.Item(44, 5).Value = 15
$ws.Cells.Item(44, 6).Value = 100112027
$ws.Cells.Item(44, 7).Value = "Melón"
$ws.Cells.Item(44, 8).Value = "Calameño"
$ws.Cells.Item(44, 9).Value = "Segunda"
$ws.Cells.Item(44, 10).Value = 100
$ws.Cells.Item(44, 11).Value = 19000
$ws.Cells.Item(44, 12).Value = 20000
$ws.Cells.Item(44, 13).Value = 19500
$ws.Cells.Item(44, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(44, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(44, 16).Value = 812
$ws.Cells.Item(44, 17).Value = 24
$ws.Cells.Item(44, 18).Value = "Hortaliza"

# Row 45: Melón, Tuna, Segunda
$ws.Cells.Item(45, 1).Value = 1
$ws.Cells.Item(45, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(45, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(45, 4).Value = 44904
$ws.Cells.Item(45, 5).Value = 15
$ws.Cells.Item(45, 6).Value = 100112027
$ws.Cells.Item(45, 7).Value = "Melón"
$ws.Cells.Item(45, 8).Value = "Tuna"
$ws.Cells.Item(45, 9).Value = "Segunda"
$ws.Cells.Item(45, 10).Value = 130
$ws.Cells.Item(45, 11).Value = 19000
$ws.Cells.Item(45, 12).Value = 20000
$ws.Cells.Item(45, 13).Value = 19500
$ws.Cells.Item(45, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(45, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(45, 16).Value = 812
$ws.Cells.Item(45, 17).Value = 24
$ws.Cells.Item(45, 18).Value = "Hortaliza"
